# feat: add 2022-Q3 data
#
# 1) Insert a new worksheet "2022-Q3" right before the current "2022-Q1"
#    sheet (all later quarter sheets shift right by one position).
# 2) Populate "2022-Q3" with the fund-holdings table, matching the
#    look (bold/centered/bordered header row + first data column) of the
#    other quarter sheets.
# 3) Update the "总计" (Total) summary sheet: insert a new row for
#    "2022-Q3" right after the header row, pushing the existing quarters
#    down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: insert the new "2022-Q3" worksheet before "2022-Q1"
# ---------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$newSheet = $wb.Worksheets.Add($q1Sheet)
$newSheet.Name = "2022-Q3"

# Match page setup used by all the other sheets in this workbook.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Header row (row 1), columns B:H - bold/centered/bordered like the
# other quarter sheets.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows 2-7.
$rows = @(
    @{ A=0; B="502023"; C="鹏华国证钢铁行业指数（LOF）A"; D="9.48"; E="94.49"; F="2.65"; G="0.2512"; H=9 },
    @{ A=1; B="012810"; C="鹏华国证钢铁行业指数（LOF）C"; D="4.34"; E="94.49"; F="2.65"; G="0.1150"; H=9 },
    @{ A=2; B="168203"; C="中融国证钢铁行业指数A";        D="3.34"; E="92.81"; F="2.60"; G="0.0868"; H=9 },
    @{ A=3; B="164811"; C="工银瑞信中证京津冀协同发展主题指数（LOF）A"; D="0.12"; E="93.09"; F="3.05"; G="0.0037"; H=6 },
    @{ A=4; B="164825"; C="工银瑞信中证京津冀协同发展主题指数（LOF）C"; D="0.03"; E="93.09"; F="3.05"; G="0.0009"; H=6 },
    @{ A=5; B="016815"; C="中融国证钢铁行业指数C";        D="0.00"; E="92.81"; F="2.60"; G=0;      H=9 }
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row.A

    # Text-style columns: force text format first so fund codes/decimal
    # strings ("502023", "9.48", ...) are not coerced into numbers and
    # don't lose leading zeros / precision, then clear the format back
    # off again (the value is already stored as text by then, so this
    # doesn't re-coerce it) so the cell ends up with the same "no
    # explicit style" look the source data uses.
    $c = $newSheet.Cells.Item($r, 2)
    $c.NumberFormat = "@"
    $c.Value = $row.B
    $c.ClearFormats()

    $c = $newSheet.Cells.Item($r, 3)
    $c.NumberFormat = "@"
    $c.Value = $row.C
    $c.ClearFormats()

    $c = $newSheet.Cells.Item($r, 4)
    $c.NumberFormat = "@"
    $c.Value = $row.D
    $c.ClearFormats()

    $c = $newSheet.Cells.Item($r, 5)
    $c.NumberFormat = "@"
    $c.Value = $row.E
    $c.ClearFormats()

    $c = $newSheet.Cells.Item($r, 6)
    $c.NumberFormat = "@"
    $c.Value = $row.F
    $c.ClearFormats()

    $c = $newSheet.Cells.Item($r, 7)
    if ($row.G -is [string]) {
        $c.NumberFormat = "@"
    }
    $c.Value = $row.G
    $c.ClearFormats()

    $newSheet.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# Copy the header/first-column formatting (bold, centered, thin border)
# from the corresponding cells on the "2022-Q1" sheet so the new sheet
# looks the same as the rest. NOTE: re-fetch "2022-Q1" by name here -
# the original $q1Sheet variable now resolves (positionally) to the
# newly-inserted sheet itself, since Add() spliced the new sheet in at
# that slot.
$q1SheetAgain = $wb.Worksheets.Item("2022-Q1")
$q1SheetAgain.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q1SheetAgain.Range("A2").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Part 2: update the "总计" (Total) sheet with the new 2022-Q3 row
# ---------------------------------------------------------------------
# Column A is just the 0-based row index (row 2 -> 0, row 3 -> 1, ...)
# and does not change; only the quarter label / count / value in
# columns B:D cascade down by one row, with a brand-new row (A=6) added
# at the bottom for "2020-Q4", which is what got pushed off the list.
$totalSheet = $wb.Worksheets.Item("总计")

# Create row 8 first by copying row 7's formatting (so A8 picks up the
# same "s=2" bordered/bold/centered style used by the rest of column A).
$totalSheet.Range("A7:D7").Copy()
$totalSheet.Range("A8:D8").PasteSpecial(-4122)
$totalSheet.Range("A8").Value = 6

$totalData = @(
    @{ B="2022-Q3"; C=6;  D=0.46 },
    @{ B="2022-Q1"; C=2;  D=0.04 },
    @{ B="2021-Q4"; C=8;  D=0.29 },
    @{ B="2021-Q3"; C=3;  D=0.96 },
    @{ B="2021-Q2"; C=5;  D=0.45 },
    @{ B="2021-Q1"; C=5;  D=0.39 },
    @{ B="2020-Q4"; C=10; D=0.74 }
)
$r = 2
foreach ($row in $totalData) {
    $totalSheet.Cells.Item($r, 2).Value = $row.B
    $totalSheet.Cells.Item($r, 3).Value = $row.C
    $totalSheet.Cells.Item($r, 4).Value = $row.D
    $r = $r + 1
}
